$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.635.32"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.75%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.852.49"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.15%  "

$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "456.89"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +7.82%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.66"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +12.37%  "

$ws.Range("E7").Value = "  +2.54%  "

$ws.Range("E8").Value = "  +0.01%  "

$ws.Range("E9").Value = "  +3.07%  "

$ws.Range("E10").Value = "  -1.75%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0000317"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -7.95%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "43.89"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +7.30%  "

$ws.Range("E13").Value = "  -0.37%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.457.05"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.22%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.77"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -6.14%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.851.57"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.19%  "

$ws.Range("E17").Value = "  -0.31%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "20.01"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.84%  "

$ws.Range("E19").Value = "  +6.77%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "67.593.08"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.35%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "427.62"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +4.33%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "14.88"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.03%  "

$ws.Range("E23").Value = "  +7.09%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "86.94"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.05%  "

$ws.Range("B25").Value = "RenderToken"
$ws.Range("C25").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "10.76"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +13.15%  "

$ws.Range("B26").Value = "PancakeSwap"
$ws.Range("C26").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.51"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +8.49%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "37.57"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.43%  "

$ws.Range("E28").Value = "  -0.42%  "

$ws.Range("E29").Value = "  +1.35%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "745.86"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.24%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.136"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +12.04%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "13.84"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +5.16%  "

$ws.Range("E33").Value = "  -1.20%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "43.46"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +12.48%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.162"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +7.30%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "57.56"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.28%  "

$ws.Range("E37").Value = "  +3.05%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.999"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.02%  "

$ws.Range("E39").Value = "  +4.56%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.359"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +13.65%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.97"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.34%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.64"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +16.12%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0₃0680"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -9.30%  "

$ws.Range("E44").Value = "  +4.88%  "

$ws.Range("E45").Value = "  +0.05%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.44"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.10%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.25"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.92%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.75"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +7.99%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.13"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.75%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "144.66"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.89%  "

$ws.Range("E51").Value = "  +2.25%  "
